$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the sheet back down to a single sample data row (the big bulk
# export used to run down to row 355997 - this is the "Removed data" /
# "Reduce test files sized" part of the commit).
$ws.Range("A2").Value = 2019

# Column widths drift very slightly (10.7109375 -> 10.6640625 etc.) as a
# side effect of the resave; nudge them to the closest reachable widths.
$ws.Range("A1:E1").EntireColumn.ColumnWidth = 9.83
$ws.Range("F1").EntireColumn.ColumnWidth = 19.83
$ws.Range("G1").EntireColumn.ColumnWidth = 9.83

# Leave the cursor where the author's save left it.
$ws.Range("A3").Select() | Out-Null
